$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 141 is the last existing data row (date serial 45697).
# Append two more rows (142, 143) that repeat the same values as row 141,
# incrementing only the date in column A by one day each time.

$srcRow = 141
$lastCol = 10  # column J

for ($i = 1; $i -le 2; $i++) {
    $newRow = $srcRow + $i

    # Copy formatting (styles) from the source row to the new row so the
    # new cells keep the same number formats/borders/fonts without
    # introducing new style entries.
    $ws.Range("A$srcRow`:J$srcRow").Copy()
    $ws.Range("A$newRow`:J$newRow").PasteSpecial(-4122)

    # Column A: increment the date serial by one day from the source row.
    $ws.Cells.Item($newRow, 1).Value = 45697 + $i

    # Columns B-J: copy the same values as the source row.
    for ($col = 2; $col -le $lastCol; $col++) {
        $ws.Cells.Item($newRow, $col).Value = $ws.Cells.Item($srcRow, $col).Value2
    }
}

$excel.CutCopyMode = 0

Write-Host "done"
